$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 508, shifting existing rows 508..549 down to 509..550
$ws.Rows.Item(508).EntireRow.Insert()

# Populate the newly inserted row 508 with the new record's data
$ws.Range("A508").Value = 5
$ws.Range("B508").Value = "Macroferia Regional de Talca"
$ws.Range("C508").Value = "Maule"
$ws.Range("D508").Value = 45106
$ws.Range("E508").Value = 7
$ws.Range("F508").Value = 100114013
$ws.Range("G508").Value = "Zanahoria"
$ws.Range("H508").Value = "Sin especificar"
$ws.Range("I508").Value = "Primera"
$ws.Range("J508").Value = 500
$ws.Range("K508").Value = 6500
$ws.Range("L508").Value = 6500
$ws.Range("M508").Value = 6500
$ws.Range("N508").Value = "$/saco 20 kilos"
$ws.Range("O508").Value = "Región de Ñuble"
$ws.Range("P508").Value = 325
$ws.Range("Q508").Value = 20
$ws.Range("R508").Value = "Hortaliza"
